$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") contains text values that look numeric (e.g. "552.34",
# "65.256.96", "0.0000261"). Force the column to Text format first so Excel
# does not silently convert the assigned strings into numbers/dates.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "65.256.96"
$ws.Range("E2").Value = "  -6.63%  "
$ws.Range("D3").Value = "3.273.30"
$ws.Range("E3").Value = "  -7.95%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "552.34"
$ws.Range("E5").Value = "  -6.58%  "
$ws.Range("D6").Value = "178.51"
$ws.Range("E6").Value = "  -9.28%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "0.584"
$ws.Range("E8").Value = "  -4.80%  "
$ws.Range("D9").Value = "3.269.91"
$ws.Range("E9").Value = "  -7.60%  "
$ws.Range("D10").Value = "0.186"
$ws.Range("E10").Value = "  -11.04%  "
$ws.Range("D11").Value = "0.582"
$ws.Range("E11").Value = "  -7.31%  "
$ws.Range("D12").Value = "46.98"
$ws.Range("E12").Value = "  -11.28%  "
$ws.Range("D13").Value = "0.0000261"
$ws.Range("E13").Value = "  -10.47%  "
$ws.Range("D14").Value = "8.48"
$ws.Range("E14").Value = "  -8.99%  "
$ws.Range("D15").Value = "3.791.66"
$ws.Range("E15").Value = "  -8.04%  "
$ws.Range("D16").Value = "604.23"
$ws.Range("E16").Value = "  -6.62%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "65.234.49"
$ws.Range("E17").Value = "  -6.64%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "17.76"
$ws.Range("E18").Value = "  -4.03%  "
$ws.Range("E19").Value = "  -4.39%  "
$ws.Range("D20").Value = "3.262.22"
$ws.Range("E20").Value = "  -8.14%  "
$ws.Range("D21").Value = "11.32"
$ws.Range("E21").Value = "  -10.21%  "
$ws.Range("D22").Value = "0.895"
$ws.Range("E22").Value = "  -7.36%  "
$ws.Range("D23").Value = "17.57"
$ws.Range("E23").Value = "  -3.46%  "
$ws.Range("D24").Value = "102.15"
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("E25").Value = "  -10.23%  "
$ws.Range("D26").Value = "3.95"
$ws.Range("E26").Value = "  -10.56%  "
$ws.Range("D27").Value = "5.98"
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("D28").Value = "2.66"
$ws.Range("E28").Value = "  -9.73%  "
$ws.Range("D29").Value = "9.32"
$ws.Range("E29").Value = "  -9.35%  "
$ws.Range("D30").Value = "8.59"
$ws.Range("E30").Value = "  -10.94%  "
$ws.Range("D31").Value = "30.25"
$ws.Range("E31").Value = "  -9.40%  "
$ws.Range("D32").Value = "3.87"
$ws.Range("E32").Value = "  -10.78%  "
$ws.Range("D33").Value = "6.20"
$ws.Range("E33").Value = "  -9.63%  "
$ws.Range("D34").Value = "10.97"
$ws.Range("E34").Value = "  -7.08%  "
$ws.Range("D35").Value = "541.33"
$ws.Range("E35").Value = "  +5.05%  "
$ws.Range("D36").Value = "0.104"
$ws.Range("E36").Value = "  -7.07%  "
$ws.Range("D37").Value = "3.745.06"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("D39").Value = "56.33"
$ws.Range("E39").Value = "  -8.89%  "
$ws.Range("D40").Value = "3.42"
$ws.Range("E40").Value = "  -6.93%  "
$ws.Range("D41").Value = "0.0₃0700"
$ws.Range("E41").Value = "  -14.07%  "
$ws.Range("D42").Value = "2.66"
$ws.Range("E42").Value = "  -10.40%  "
$ws.Range("D43").Value = "0.125"
$ws.Range("E43").Value = "  -6.95%  "
$ws.Range("D44").Value = "31.64"
$ws.Range("E44").Value = "  -10.04%  "
$ws.Range("B45").Value = "CoreDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D45").Value = "3.23"
$ws.Range("E45").Value = "  +18.10%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "0.335"
$ws.Range("E46").Value = "  -10.27%  "
$ws.Range("D47").Value = "0.0407"
$ws.Range("E47").Value = "  -10.59%  "
$ws.Range("D48").Value = "3.19"
$ws.Range("E48").Value = "  -7.05%  "
$ws.Range("E49").Value = "  -6.59%  "
$ws.Range("D50").Value = "2.57"
$ws.Range("E50").Value = "  -11.20%  "
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").Value = "  -0.17%  "
